# Update crypto price/volume data (scraped values refreshed by GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so decimal-looking price strings
# (e.g. "1.002") are not silently reinterpreted as numbers, then clear the
# temporary formatting again once all values are written so no residual
# style/number-format is left on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.868.94'
$ws.Range("E2").Value = '  -2.66%  '

$ws.Range("D3").Value = '1.886.57'

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -1.08%  '

$ws.Range("E5").Value = '  +1.70%  '

$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.92%  '

$ws.Range("D7").Value = '0.4582'
$ws.Range("E7").Value = '  -3.98%  '

$ws.Range("D8").Value = '0.3925'
$ws.Range("E8").Value = '  -2.99%  '

$ws.Range("D9").Value = '49.31'
$ws.Range("E9").Value = '  -8.42%  '

$ws.Range("D10").Value = '0.08242'
$ws.Range("E10").Value = '  -2.74%  '

$ws.Range("D11").Value = '1.035'
$ws.Range("E11").Value = '  -2.44%  '

$ws.Range("D12").Value = '21.85'
$ws.Range("E12").Value = '  -1.83%  '

$ws.Range("D13").Value = '1.889.96'
$ws.Range("E13").Value = '  -2.57%  '

$ws.Range("D14").Value = '7.309'
$ws.Range("E14").Value = '  -4.14%  '

$ws.Range("D15").Value = '5.964'
$ws.Range("E15").Value = '  -4.05%  '

$ws.Range("E16").Value = '  -0.95%  '

$ws.Range("D17").Value = '88.83'
$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("E18").Value = '  -1.95%  '

$ws.Range("D19").Value = '0.06574'
$ws.Range("E19").Value = '  -0.77%  '

$ws.Range("D20").Value = '17.38'
$ws.Range("E20").Value = '  -7.03%  '

$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  -0.92%  '

$ws.Range("D22").Value = '5.615'
$ws.Range("E22").Value = '  -3.77%  '

$ws.Range("D23").Value = '27.896.87'

$ws.Range("D24").Value = '11.05'
$ws.Range("E24").Value = '  -4.29%  '

$ws.Range("D25").Value = '2.301'
$ws.Range("E25").Value = '  +0.22%  '

$ws.Range("D26").Value = '2.154.31'
$ws.Range("E26").Value = '  -0.81%  '

$ws.Range("D27").Value = '153.96'
$ws.Range("E27").Value = '  -0.34%  '

$ws.Range("E28").Value = '  -1.86%  '

$ws.Range("D29").Value = '5.695'
$ws.Range("E29").Value = '  -4.39%  '

$ws.Range("D30").Value = '2.093'
$ws.Range("E30").Value = '  -3.01%  '

$ws.Range("D31").Value = '123.21'
$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("D32").Value = '0.09531'
$ws.Range("E32").Value = '  -0.58%  '

$ws.Range("D33").Value = '0.9547'
$ws.Range("E33").Value = '  -4.84%  '

$ws.Range("D34").Value = '1.474'
$ws.Range("E34").Value = '  +2.33%  '

$ws.Range("D35").Value = '3.633'
$ws.Range("E35").Value = '  -0.85%  '

$ws.Range("D36").Value = '5.432'
$ws.Range("E36").Value = '  -4.05%  '

$ws.Range("D37").Value = '0.02275'
$ws.Range("E37").Value = '  -2.96%  '

$ws.Range("D38").Value = '1.249'
$ws.Range("E38").Value = '  -1.71%  '

$ws.Range("D39").Value = '0.06087'
$ws.Range("E39").Value = '  -2.35%  '

$ws.Range("D40").Value = '8.567'
$ws.Range("E40").Value = '  -1.81%  '

$ws.Range("D41").Value = '0.6091'
$ws.Range("E41").Value = '  -2.39%  '

$ws.Range("E42").Value = '  -0.97%  '

$ws.Range("D43").Value = '10.71'
$ws.Range("E43").Value = '  -3.64%  '

$ws.Range("D44").Value = '0.1897'
$ws.Range("E44").Value = '  -1.29%  '

$ws.Range("D45").Value = '1.301'
$ws.Range("E45").Value = '  -3.14%  '

$ws.Range("E46").Value = '  -2.43%  '

$ws.Range("D47").Value = '12.69'
$ws.Range("E47").Value = '  -1.89%  '

$ws.Range("D48").Value = '1.985'
$ws.Range("E48").Value = '  -4.72%  '

$ws.Range("E49").Value = '  +0.48%  '

$ws.Range("D50").Value = '0.06903'
$ws.Range("E50").Value = '  +1.22%  '

$ws.Range("D51").Value = '110.08'
$ws.Range("E51").Value = '  -1.12%  '

$priceRange.ClearFormats()
